$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new Job Posting row (row 10): Job_Id=9
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "QA Analyst"
$ws.Cells.Item(10, 3).Value = "Testing1"
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0
